$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Mapping sheet: shapefiles were re-projected to WGS 84 (EPSG:4269), so the
# reserve / station bounding-box coordinates change. Previously A and B held
# identical values (a single bound); now they hold distinct left/right
# bounds.
# ---------------------------------------------------------------------------
$mapping = $wb.Worksheets.Item("Mapping")

$mapping.Range("A2").Value = -81.4486
$mapping.Range("B2").Value = -81.4175

$mapping.Range("A3").Value = 29.5473
$mapping.Range("B3").Value = 29.5744

$mapping.Range("A4").Value = -81.1353
$mapping.Range("B4").Value = -81.1663

$mapping.Range("A5").Value = 30.2152
$mapping.Range("B5").Value = 30.1882

# ---------------------------------------------------------------------------
# Basic_Plotting sheet: "atemp" (row 8) had its Parameter_Category fixed up.
# ---------------------------------------------------------------------------
$basicPlotting = $wb.Worksheets.Item("Basic_Plotting")
$basicPlotting.Range("B8").Value = "instantaneous"
